$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H7").Value = 1900
$ws.Range("I7").Value = 666.6667
$ws.Range("J7").Value = 3750
$ws.Range("K7").Value = 666.6667
$ws.Range("L7").Value = 3750
$ws.Range("M7").Value = -554.6667
$ws.Range("N7").Value = -3974
$ws.Range("H8").Value = 34.2
$ws.Range("I8").Value = 34.2
$ws.Range("K8").Value = 102.6
$ws.Range("M8").Value = 36.39999999999999
$ws.Range("H12").Value = 421
$ws.Range("J12").Value = 500.66666
$ws.Range("L12").Value = 500.66666
$ws.Range("N12").Value = -840.66666
$ws.Range("H14").Value = 1900
$ws.Range("I14").Value = 666.6667
$ws.Range("J14").Value = 3750
$ws.Range("K14").Value = 666.6667
$ws.Range("L14").Value = 3750
$ws.Range("M14").Value = -475.6667
$ws.Range("N14").Value = -4132
$ws.Range("H62").Value = 6846.6665
$ws.Range("J62").Value = 6995
$ws.Range("L62").Value = 6995
$ws.Range("N62").Value = -8243
$ws.Range("H65").Value = 6846.6665
$ws.Range("J65").Value = 6995
$ws.Range("L65").Value = 34975
$ws.Range("N65").Value = -41215
$ws.Range("H70").Value = 15499.667
$ws.Range("J70").Value = 18499.572
$ws.Range("L70").Value = 55498.716
$ws.Range("N70").Value = -56038.716
$ws.Range("H73").Value = 15499.667
$ws.Range("J73").Value = 18499.572
$ws.Range("L73").Value = 55498.716
$ws.Range("N73").Value = -57370.716
$ws.Range("H74").Value = 8376
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 8376
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 8376
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -10248
$ws.Range("H77").Value = 8376
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 8376
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 41880
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -51240
$ws.Range("H132").Value = 2064.5386
$ws.Range("I132").Value = 2064.5386
$ws.Range("K132").Value = 6193.6158
$ws.Range("M132").Value = -3663.6158
$ws.Range("H138").Value = 4157.0835
$ws.Range("J138").Value = 5569.2856
$ws.Range("L138").Value = 16707.8568
$ws.Range("N138").Value = -26987.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 29233.2
$ws.Range("J101").Value = 29233.2
$ws.Range("L101").Value = 29233.2
$ws.Range("N101").Value = -35723.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 40591.4
$ws.Range("I86").Value = 1475
$ws.Range("K86").Value = 1475
$ws.Range("M86").Value = -352
$ws.Range("H89").Value = 40591.4
$ws.Range("I89").Value = 1475
$ws.Range("K89").Value = 7375
$ws.Range("M89").Value = -1759
$ws.Range("H99").Value = 3406.6
$ws.Range("I99").Value = 3008.5
$ws.Range("K99").Value = 3008.5
$ws.Range("M99").Value = -1510.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1102
$ws.Range("I17").Value = 802.6667
$ws.Range("K17").Value = 802.6667
$ws.Range("M17").Value = -628.6667
$ws.Range("H31").Value = 2295.6924
$ws.Range("J31").Value = 3997.4
$ws.Range("L31").Value = 3997.4
$ws.Range("N31").Value = -4587.4
$ws.Range("H34").Value = 2295.6924
$ws.Range("J34").Value = 3997.4
$ws.Range("L34").Value = 3997.4
$ws.Range("N34").Value = -4401.4
$ws.Range("H107").Value = 785.5294
$ws.Range("I107").Value = 381.58334
$ws.Range("J107").Value = 1755
$ws.Range("K107").Value = 381.58334
$ws.Range("L107").Value = 1755
$ws.Range("M107").Value = 1538.41666
$ws.Range("N107").Value = -5595
$ws.Range("H122").Value = 1191.5454
$ws.Range("I122").Value = 1149.3334
$ws.Range("J122").Value = 1242.2
$ws.Range("K122").Value = 3448.0002
$ws.Range("L122").Value = 3726.6
$ws.Range("M122").Value = -998.0001999999999
$ws.Range("N122").Value = -8626.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 384.125
$ws.Range("J12").Value = 17.75
$ws.Range("L12").Value = 53.25
$ws.Range("N12").Value = -399.25
$ws.Range("H131").Value = 1858.1666
$ws.Range("I131").Value = 849.75
$ws.Range("J131").Value = 2362.375
$ws.Range("K131").Value = 2549.25
$ws.Range("L131").Value = 7087.125
$ws.Range("M131").Value = 2490.75
$ws.Range("N131").Value = -17167.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 286235.84
$ws.Range("I3").Value = 500640.6
$ws.Range("J3").Value = 167122.11
$ws.Range("K3").Value = 500640.6
$ws.Range("L3").Value = 167122.11
$ws.Range("M3").Value = -500524.6
$ws.Range("N3").Value = -167354.11
$ws.Range("H26").Value = 30000
$ws.Range("J26").Value = 30000
$ws.Range("L26").Value = 30000
$ws.Range("N26").Value = -30560
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -30996
$ws.Range("H80").Value = 3074.3333
$ws.Range("I80").Value = 3074.3333
$ws.Range("K80").Value = 3074.3333
$ws.Range("M80").Value = -2076.3333
$ws.Range("H83").Value = 3074.3333
$ws.Range("I83").Value = 3074.3333
$ws.Range("K83").Value = 15371.6665
$ws.Range("M83").Value = -10379.6665
$ws.Range("H101").Value = 44666
$ws.Range("J101").Value = 44666
$ws.Range("L101").Value = 44666
$ws.Range("N101").Value = -51156
$ws.Range("H126").Value = 4316.1
$ws.Range("I126").Value = 2753.6667
$ws.Range("K126").Value = 8261.000100000001
$ws.Range("M126").Value = -5791.000100000001
$ws.Range("H132").Value = 6859.5713
$ws.Range("I132").Value = 5604.6
$ws.Range("K132").Value = 16813.8
$ws.Range("M132").Value = -14283.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4524.278
$ws.Range("I7").Value = 4477.3125
$ws.Range("J7").Value = 4900
$ws.Range("K7").Value = 4477.3125
$ws.Range("L7").Value = 4900
$ws.Range("M7").Value = -4365.3125
$ws.Range("N7").Value = -5124
$ws.Range("H68").Value = 4875
$ws.Range("I68").Value = 2714.2856
$ws.Range("J68").Value = 20000
$ws.Range("K68").Value = 2714.2856
$ws.Range("L68").Value = 20000
$ws.Range("M68").Value = -1965.2856
$ws.Range("N68").Value = -21498
$ws.Range("H71").Value = 4875
$ws.Range("I71").Value = 2714.2856
$ws.Range("J71").Value = 20000
$ws.Range("K71").Value = 13571.428
$ws.Range("L71").Value = 100000
$ws.Range("M71").Value = -9827.428
$ws.Range("N71").Value = -107488
$ws.Range("H82").Value = 2457.3845
$ws.Range("I82").Value = 2118.5
$ws.Range("J82").Value = 2999.6
$ws.Range("K82").Value = 2118.5
$ws.Range("L82").Value = 2999.6
$ws.Range("M82").Value = -1757.5
$ws.Range("N82").Value = -3721.6
$ws.Range("H85").Value = 2457.3845
$ws.Range("I85").Value = 2118.5
$ws.Range("J85").Value = 2999.6
$ws.Range("K85").Value = 2118.5
$ws.Range("L85").Value = 2999.6
$ws.Range("M85").Value = -870.5
$ws.Range("N85").Value = -5495.6
$ws.Range("H126").Value = 4524.278
$ws.Range("I126").Value = 4477.3125
$ws.Range("J126").Value = 4900
$ws.Range("K126").Value = 13431.9375
$ws.Range("L126").Value = 14700
$ws.Range("M126").Value = -10961.9375
$ws.Range("N126").Value = -19640
$ws.Range("H132").Value = 9699.799999999999
$ws.Range("I132").Value = 8874.75
$ws.Range("J132").Value = 13000
$ws.Range("K132").Value = 26624.25
$ws.Range("L132").Value = 39000
$ws.Range("M132").Value = -24094.25
$ws.Range("N132").Value = -44060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 2733.1667
$ws.Range("I132").Value = 2939.8
$ws.Range("J132").Value = 1700
$ws.Range("K132").Value = 8819.400000000001
$ws.Range("L132").Value = 5100
$ws.Range("M132").Value = -6289.400000000001
$ws.Range("N132").Value = -10160
